# daily auto push: 2025-10-01 09:30 UTC
# Append the next daily data row to Sheet1 (row 45):
#   A45 = "2025/10/01"  (date recorded as plain text, same as column A elsewhere)
#   B45 = "水"          (weekday label, text)
#   C45 = 17            (hour, number)
#   D45 = 17            (ranking, number)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = $ws.Cells.Item($ws.UsedRange.Rows.Count, 1).Row + 1

# Column A holds the date as literal text (e.g. "2025/09/22") rather than a
# real date value elsewhere in the sheet, so force text entry here too.
# Applying a temporary "@" (text) number format stops Excel's automatic
# date-literal recognition from converting the string into a date serial;
# resetting the cell style back to "Normal" afterwards drops the temporary
# formatting so the new cell ends up unstyled, matching the rest of the
# column.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/01"
$cellA.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "水"
$ws.Cells.Item($newRow, 3).Value = 17
$ws.Cells.Item($newRow, 4).Value = 17
